$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timesheet")

# --- Update time entries (values in minutes column F, hours column D) ---
# Row 9: "Camara perspectiva fixa"
$ws.Range("F9").Value = 30

# Row 10: "Camara perspectiva movel"
$ws.Range("F10").Value = 30

# Row 11: "Seleccao de camara activa por teclado"
$ws.Range("F11").Value = 30

# Row 12: "Movimento das laranjas"
$ws.Range("D12").Value = 15

# Row 13: "Aparecimento / Desaparecimento das Laranjas"
$ws.Range("F13").Value = 30

# Row 14: "Colisoes entre carro e laranjas + perda de vida"
$ws.Range("F14").Value = 30

# Row 15: "Colisoes entre carro e cheerios e entre cheerios"
$ws.Range("D15").Value = 26

# Row 16: "Movimento dos Cheerios"
$ws.Range("F16").Value = 5

# --- Update the selected / visible cell in the sheet view ---
$ws.Range("D16").Select()

$wb.Save()
